# Apply updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.438.99'
$ws.Range('E2').Value = '  +0.23%  '

$ws.Range('D3').Value = '1.573.73'
$ws.Range('E3').Value = '  +0.12%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('E5').Value = '  -0.03%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.18%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3732'
$ws.Range('D7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.98'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.13%  '

$ws.Range('E9').Value = '  -0.83%  '

$ws.Range('E10').Value = '  -0.76%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.145'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.65%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.11%  '

$ws.Range('E13').Value = '  +0.73%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.021'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.12%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.966'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.38%  '

$ws.Range('D16').Value = '1.573.87'
$ws.Range('E16').Value = '  +0.20%  '

$ws.Range('E17').Value = '  -0.66%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.96%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06752'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.08%  '

$ws.Range('E20').Value = '  +0.06%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.309'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.80%  '

$ws.Range('E22').Value = '  -2.90%  '

$ws.Range('E23').Value = '  +1.36%  '

$ws.Range('D24').Value = '22.431.06'
$ws.Range('E24').Value = '  +0.16%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.343'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.25%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.688'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.34%  '

$ws.Range('E27').Value = '  -0.30%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '148.61'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.88%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.009'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.40%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.41%  '

$ws.Range('D31').Value = '1.749.45'
$ws.Range('E31').Value = '  +0.15%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.057'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.41%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.219'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.26%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.984'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.11%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.841'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.84%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08385'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.83%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.374'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.57%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02491'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.13%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2295'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.92%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06525'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.54%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.491'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.43%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.32'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.94%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6239'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.45%  '

$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.01%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.26%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.812'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.76%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5826'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.60%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.54'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.53%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.080'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.27%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.228'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.62%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07330'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.07%  '
